$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1169995834814548
$ws.Range("C2").Value = 109.9114832445916
$ws.Range("D2").Value = 189.6080260415259
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 313.5003553403997

# Row 3
$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3993.344853322108
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 4009.490252764333

# Row 4
$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 4.327115817150455
